# Apply the "Şube-Öğretmen" (Branch-Teacher) table addition to Sayfa1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")  # "Sayfa1" is the active/tabSelected sheet

# Update existing header labels to the more specific "*Id" wording.
$ws.Range("H2").Value = "Şube Id"
$ws.Range("G9").Value = "Şube Id"
$ws.Range("H9").Value = "Branş Id"

# New "Şube-Öğretmen" mapping table, to the right of the "Şube" table (rows 15-20).
# Match styling of the neighboring "Şube" table title/header cells first (Copy +
# PasteSpecial formats reuses the existing style record instead of minting a new one).
$ws.Range("A15").Copy()
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("A16").Copy()
$ws.Range("E16:G16").PasteSpecial(-4122)

$ws.Range("E15").Value = "Şube-Öğretmen"

$ws.Range("E16").Value = "Şube Id"
$ws.Range("F16").Value = "Öğretmen Id"
$ws.Range("G16").Value = "Ders Id"

$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1

$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 1

$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 3
$ws.Range("G19").Value = 2

$ws.Range("E20").Value = 4
$ws.Range("F20").Value = 4

# Update selection state to match the saved view (also clears the old
# "topLeftCell=A10" scroll position, matching the diff).
$ws.Range("G24").Select()
